$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitor Order")
$ws.Activate()
$win = $app.ActiveWindow
try {
  $win.TopLeftCell = $ws.Range("A4")
  Write-Output "topleftcell set ok"
} catch {
  Write-Output "err1: $_"
}
try {
  $ws.Range("A4").TopLeftCell.Select()
} catch {
  Write-Output "err2: $_"
}
